# Applies numeric updates to the per-job "Leve Profits" tables across all eight
# crafting-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), refreshing
# cached market-board price/profit figures pulled in by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Cells.Item(41, 8).Value = 352.9565  # H41
$ws.Cells.Item(41, 9).Value = 344.8  # I41
$ws.Cells.Item(41, 11).Value = 344.8  # K41
$ws.Cells.Item(41, 13).Value = 95.19999999999999  # M41
$ws.Cells.Item(62, 8).Value = 16675892  # H62
$ws.Cells.Item(62, 9).Value = 23820564  # I62
$ws.Cells.Item(62, 11).Value = 23820564  # K62
$ws.Cells.Item(62, 13).Value = -23819940  # M62
$ws.Cells.Item(65, 8).Value = 16675892  # H65
$ws.Cells.Item(65, 9).Value = 23820564  # I65
$ws.Cells.Item(65, 11).Value = 119102820  # K65
$ws.Cells.Item(65, 13).Value = -119099700  # M65
$ws.Cells.Item(94, 8).Value = 4640.25  # H94
$ws.Cells.Item(94, 9).Value = 874.7143  # I94
$ws.Cells.Item(94, 11).Value = 874.7143  # K94
$ws.Cells.Item(94, 13).Value = -423.7143  # M94
$ws.Cells.Item(124, 8).Value = 250000  # H124
$ws.Cells.Item(124, 9).Value = 100000  # I124
$ws.Cells.Item(124, 10).Value = 400000  # J124
$ws.Cells.Item(124, 11).Value = 100000  # K124
$ws.Cells.Item(124, 12).Value = 400000  # L124
$ws.Cells.Item(124, 13).Value = -95090  # M124
$ws.Cells.Item(124, 14).Value = -409820  # N124
$ws.Cells.Item(132, 8).Value = 226890.92  # H132
$ws.Cells.Item(132, 9).Value = 259861.64  # I132
$ws.Cells.Item(132, 11).Value = 779584.92  # K132
$ws.Cells.Item(132, 13).Value = -777054.92  # M132
$ws.Cells.Item(137, 8).Value = 3945.2334  # H137
$ws.Cells.Item(137, 9).Value = 3721.4211  # I137
$ws.Cells.Item(137, 10).Value = 4331.8184  # J137
$ws.Cells.Item(137, 11).Value = 11164.2633  # K137
$ws.Cells.Item(137, 12).Value = 12995.4552  # L137
$ws.Cells.Item(137, 13).Value = -8614.263300000001  # M137
$ws.Cells.Item(137, 14).Value = -18095.4552  # N137
$ws.Cells.Item(139, 8).Value = 157778.92  # H139
$ws.Cells.Item(139, 10).Value = 157778.92  # J139
$ws.Cells.Item(139, 12).Value = 157778.92  # L139
$ws.Cells.Item(139, 14).Value = -168058.92  # N139
$ws.Cells.Item(141, 8).Value = 3284.5  # H141
$ws.Cells.Item(141, 9).Value = 3011  # I141
$ws.Cells.Item(141, 10).Value = 4105  # J141
$ws.Cells.Item(141, 11).Value = 9033  # K141
$ws.Cells.Item(141, 12).Value = 12315  # L141
$ws.Cells.Item(141, 13).Value = -3853  # M141
$ws.Cells.Item(141, 14).Value = -22675  # N141

$ws = $wb.Worksheets.Item("ARM")

$ws.Cells.Item(32, 8).Value = 1688218.6  # H32
$ws.Cells.Item(32, 9).Value = 3793.6724  # I32
$ws.Cells.Item(32, 11).Value = 3793.6724  # K32
$ws.Cells.Item(32, 13).Value = -3506.6724  # M32
$ws.Cells.Item(61, 8).Value = 3187.82  # H61
$ws.Cells.Item(61, 9).Value = 2227.0715  # I61
$ws.Cells.Item(61, 10).Value = 4410.591  # J61
$ws.Cells.Item(61, 11).Value = 2227.0715  # K61
$ws.Cells.Item(61, 12).Value = 4410.591  # L61
$ws.Cells.Item(61, 13).Value = -2015.0715  # M61
$ws.Cells.Item(61, 14).Value = -4834.591  # N61
$ws.Cells.Item(122, 8).Value = 6161  # H122
$ws.Cells.Item(122, 9).Value = 5098.125  # I122
$ws.Cells.Item(122, 10).Value = 7578.1665  # J122
$ws.Cells.Item(122, 11).Value = 15294.375  # K122
$ws.Cells.Item(122, 12).Value = 22734.4995  # L122
$ws.Cells.Item(122, 13).Value = -12844.375  # M122
$ws.Cells.Item(122, 14).Value = -27634.4995  # N122
$ws.Cells.Item(132, 8).Value = 687234.1  # H132
$ws.Cells.Item(132, 9).Value = 803506.4399999999  # I132
$ws.Cells.Item(132, 11).Value = 2410519.32  # K132
$ws.Cells.Item(132, 13).Value = -2407989.32  # M132
$ws.Cells.Item(136, 8).Value = 3187.82  # H136
$ws.Cells.Item(136, 9).Value = 2227.0715  # I136
$ws.Cells.Item(136, 10).Value = 4410.591  # J136
$ws.Cells.Item(136, 11).Value = 6681.2145  # K136
$ws.Cells.Item(136, 12).Value = 13231.773  # L136
$ws.Cells.Item(136, 13).Value = -4131.2145  # M136
$ws.Cells.Item(136, 14).Value = -18331.773  # N136

$ws = $wb.Worksheets.Item("BSM")

$ws.Cells.Item(134, 8).Value = 954491.4  # H134
$ws.Cells.Item(134, 9).Value = 1094355.1  # I134
$ws.Cells.Item(134, 10).Value = 10411.375  # J134
$ws.Cells.Item(134, 11).Value = 3283065.3  # K134
$ws.Cells.Item(134, 12).Value = 31234.125  # L134
$ws.Cells.Item(134, 13).Value = -3280530.3  # M134
$ws.Cells.Item(134, 14).Value = -36304.125  # N134

$ws = $wb.Worksheets.Item("CRP")

$ws.Cells.Item(16, 8).Value = 31253790  # H16
$ws.Cells.Item(16, 9).Value = 35716404  # I16
$ws.Cells.Item(16, 11).Value = 35716404  # K16
$ws.Cells.Item(16, 13).Value = -35716117  # M16
$ws.Cells.Item(31, 8).Value = 7625.5835  # H31
$ws.Cells.Item(31, 9).Value = 17081.416  # I31
$ws.Cells.Item(31, 11).Value = 17081.416  # K31
$ws.Cells.Item(31, 13).Value = -16786.416  # M31
$ws.Cells.Item(34, 8).Value = 7625.5835  # H34
$ws.Cells.Item(34, 9).Value = 17081.416  # I34
$ws.Cells.Item(34, 11).Value = 17081.416  # K34
$ws.Cells.Item(34, 13).Value = -16879.416  # M34
$ws.Cells.Item(99, 8).Value = 6949070.5  # H99
$ws.Cells.Item(99, 9).Value = 11115562  # I99
$ws.Cells.Item(99, 10).Value = 4918.5  # J99
$ws.Cells.Item(99, 11).Value = 11115562  # K99
$ws.Cells.Item(99, 12).Value = 4918.5  # L99
$ws.Cells.Item(99, 13).Value = -11114064  # M99
$ws.Cells.Item(99, 14).Value = -7914.5  # N99
$ws.Cells.Item(107, 8).Value = 329.17648  # H107
$ws.Cells.Item(107, 9).Value = 265.6  # I107
$ws.Cells.Item(107, 10).Value = 806  # J107
$ws.Cells.Item(107, 11).Value = 265.6  # K107
$ws.Cells.Item(107, 12).Value = 806  # L107
$ws.Cells.Item(107, 13).Value = 1654.4  # M107
$ws.Cells.Item(107, 14).Value = -4646  # N107
$ws.Cells.Item(113, 8).Value = 31253790  # H113
$ws.Cells.Item(113, 9).Value = 35716404  # I113
$ws.Cells.Item(113, 11).Value = 35716404  # K113
$ws.Cells.Item(113, 13).Value = -35714234  # M113
$ws.Cells.Item(126, 8).Value = 6949070.5  # H126
$ws.Cells.Item(126, 9).Value = 11115562  # I126
$ws.Cells.Item(126, 10).Value = 4918.5  # J126
$ws.Cells.Item(126, 11).Value = 33346686  # K126
$ws.Cells.Item(126, 12).Value = 14755.5  # L126
$ws.Cells.Item(126, 13).Value = -33344216  # M126
$ws.Cells.Item(126, 14).Value = -19695.5  # N126
$ws.Cells.Item(134, 8).Value = 9994.208000000001  # H134
$ws.Cells.Item(134, 9).Value = 4421.467  # I134
$ws.Cells.Item(134, 10).Value = 19282.111  # J134
$ws.Cells.Item(134, 11).Value = 13264.401  # K134
$ws.Cells.Item(134, 12).Value = 57846.333  # L134
$ws.Cells.Item(134, 13).Value = -10729.401  # M134
$ws.Cells.Item(134, 14).Value = -62916.333  # N134
$ws.Cells.Item(135, 8).Value = 0  # H135
$ws.Cells.Item(135, 10).Value = 0  # J135
$ws.Cells.Item(135, 12).Value = 0  # L135
$ws.Cells.Item(135, 14).ClearContents()  # N135

$ws = $wb.Worksheets.Item("CUL")

$ws.Cells.Item(29, 8).Value = 1654.7  # H29
$ws.Cells.Item(29, 9).Value = 2018.625  # I29
$ws.Cells.Item(29, 10).Value = 199  # J29
$ws.Cells.Item(29, 11).Value = 6055.875  # K29
$ws.Cells.Item(29, 12).Value = 597  # L29
$ws.Cells.Item(29, 13).Value = -5778.875  # M29
$ws.Cells.Item(29, 14).Value = -1151  # N29

$ws = $wb.Worksheets.Item("GSM")

$ws.Cells.Item(31, 8).Value = 2709.2856  # H31
$ws.Cells.Item(31, 9).Value = 1505  # I31
$ws.Cells.Item(31, 10).Value = 9935  # J31
$ws.Cells.Item(31, 11).Value = 1505  # K31
$ws.Cells.Item(31, 12).Value = 9935  # L31
$ws.Cells.Item(31, 13).Value = -1213  # M31
$ws.Cells.Item(31, 14).Value = -10519  # N31
$ws.Cells.Item(37, 8).Value = 2709.2856  # H37
$ws.Cells.Item(37, 9).Value = 1505  # I37
$ws.Cells.Item(37, 10).Value = 9935  # J37
$ws.Cells.Item(37, 11).Value = 1505  # K37
$ws.Cells.Item(37, 12).Value = 9935  # L37
$ws.Cells.Item(37, 13).Value = -1228  # M37
$ws.Cells.Item(37, 14).Value = -10489  # N37
$ws.Cells.Item(102, 8).Value = 4714.483  # H102
$ws.Cells.Item(102, 9).Value = 3404.2368  # I102
$ws.Cells.Item(102, 11).Value = 3404.2368  # K102
$ws.Cells.Item(102, 13).Value = -1782.2368  # M102
$ws.Cells.Item(126, 8).Value = 15160986  # H126
$ws.Cells.Item(126, 9).Value = 23815260  # I126
$ws.Cells.Item(126, 10).Value = 16008.167  # J126
$ws.Cells.Item(126, 11).Value = 71445780  # K126
$ws.Cells.Item(126, 12).Value = 48024.501  # L126
$ws.Cells.Item(126, 13).Value = -71443310  # M126
$ws.Cells.Item(126, 14).Value = -52964.501  # N126
$ws.Cells.Item(132, 8).Value = 4622.1875  # H132
$ws.Cells.Item(132, 9).Value = 4719.3  # I132
$ws.Cells.Item(132, 10).Value = 4136.625  # J132
$ws.Cells.Item(132, 11).Value = 14157.9  # K132
$ws.Cells.Item(132, 12).Value = 12409.875  # L132
$ws.Cells.Item(132, 13).Value = -11627.9  # M132
$ws.Cells.Item(132, 14).Value = -17469.875  # N132

$ws = $wb.Worksheets.Item("LTW")

$ws.Cells.Item(7, 8).Value = 8075.8237  # H7
$ws.Cells.Item(7, 9).Value = 5450.143  # I7
$ws.Cells.Item(7, 10).Value = 20329  # J7
$ws.Cells.Item(7, 11).Value = 5450.143  # K7
$ws.Cells.Item(7, 12).Value = 20329  # L7
$ws.Cells.Item(7, 13).Value = -5338.143  # M7
$ws.Cells.Item(7, 14).Value = -20553  # N7
$ws.Cells.Item(55, 8).Value = 4991.5884  # H55
$ws.Cells.Item(55, 9).Value = 3287.2727  # I55
$ws.Cells.Item(55, 10).Value = 8116.1665  # J55
$ws.Cells.Item(55, 11).Value = 3287.2727  # K55
$ws.Cells.Item(55, 12).Value = 8116.1665  # L55
$ws.Cells.Item(55, 13).Value = -3114.2727  # M55
$ws.Cells.Item(55, 14).Value = -8462.166499999999  # N55
$ws.Cells.Item(122, 8).Value = 1818471.9  # H122
$ws.Cells.Item(122, 9).Value = 2496773.5  # I122
$ws.Cells.Item(122, 10).Value = 9667  # J122
$ws.Cells.Item(122, 11).Value = 7490320.5  # K122
$ws.Cells.Item(122, 12).Value = 29001  # L122
$ws.Cells.Item(122, 13).Value = -7487870.5  # M122
$ws.Cells.Item(122, 14).Value = -33901  # N122
$ws.Cells.Item(126, 8).Value = 8075.8237  # H126
$ws.Cells.Item(126, 9).Value = 5450.143  # I126
$ws.Cells.Item(126, 10).Value = 20329  # J126
$ws.Cells.Item(126, 11).Value = 16350.429  # K126
$ws.Cells.Item(126, 12).Value = 60987  # L126
$ws.Cells.Item(126, 13).Value = -13880.429  # M126
$ws.Cells.Item(126, 14).Value = -65927  # N126
$ws.Cells.Item(132, 8).Value = 3927.111  # H132
$ws.Cells.Item(132, 9).Value = 3927.111  # I132
$ws.Cells.Item(132, 10).Value = 0  # J132
$ws.Cells.Item(132, 11).Value = 11781.333  # K132
$ws.Cells.Item(132, 12).Value = 0  # L132
$ws.Cells.Item(132, 13).Value = -9251.332999999999  # M132
$ws.Cells.Item(132, 14).ClearContents()  # N132

$ws = $wb.Worksheets.Item("WVR")

$ws.Cells.Item(107, 8).Value = 25001294  # H107
$ws.Cells.Item(107, 10).Value = 1563  # J107
$ws.Cells.Item(107, 12).Value = 4689  # L107
$ws.Cells.Item(107, 14).Value = -8529  # N107
$ws.Cells.Item(122, 8).Value = 9597.037  # H122
$ws.Cells.Item(122, 9).Value = 3209.7646  # I122
$ws.Cells.Item(122, 11).Value = 9629.293799999999  # K122
$ws.Cells.Item(122, 13).Value = -7179.293799999999  # M122
$ws.Cells.Item(126, 8).Value = 8374.666999999999  # H126
$ws.Cells.Item(126, 9).Value = 2500  # I126
$ws.Cells.Item(126, 10).Value = 9549.6  # J126
$ws.Cells.Item(126, 11).Value = 7500  # K126
$ws.Cells.Item(126, 12).Value = 28648.8  # L126
$ws.Cells.Item(126, 13).Value = -5030  # M126
$ws.Cells.Item(126, 14).Value = -33588.8  # N126
$ws.Cells.Item(132, 8).Value = 6996.7173  # H132
$ws.Cells.Item(132, 9).Value = 7445.304  # I132
$ws.Cells.Item(132, 10).Value = 6548.1304  # J132
$ws.Cells.Item(132, 11).Value = 22335.912  # K132
$ws.Cells.Item(132, 12).Value = 19644.3912  # L132
$ws.Cells.Item(132, 13).Value = -19805.912  # M132
$ws.Cells.Item(132, 14).Value = -24704.3912  # N132
$ws.Cells.Item(136, 8).Value = 10425894  # H136
$ws.Cells.Item(136, 9).Value = 16138616  # I136
$ws.Cells.Item(136, 11).Value = 48415848  # K136
$ws.Cells.Item(136, 13).Value = -48413298  # M136
